# Auto - Update data with bot!
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("D4").Value = "무료 E-Book 소개 - Python, 데이터분석, SQL, 정규식 표현"
$ws.Range("E4").Value = "https://teddylee777.github.io/thoughts/free-ebook"

# Row 16
$ws.Range("D16").Value = "Group-CAM: Group Score-Weighted Visual Explanations for Deep Convolutional Networks 내용 정리 [XAI-11]"
$ws.Range("E16").Value = "https://wewinserv.tistory.com/153"

# Row 28
$ws.Range("D28").Value = "강화학습 팁 모음"
$ws.Range("E28").Value = "https://ropiens.tistory.com/132"

# Row 32
$ws.Range("D32").Value = "배깅과 페이스팅 (Bagging, pasting)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/323"

# Row 35
$ws.Range("D35").Value = "The missing paper"

# Row 51
$ws.Range("D51").Value = "[pyqt5] 프로그램창을 항상 가장 위에 있게 하면서 동시에 타이틀 바도 없게 하려면?"
$ws.Range("E51").Value = "https://bskyvision.com/1209"
